$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2, shifting the existing
# rows 2-5 down to rows 3-6 (dimension grows from A1:N5 to A1:N6).
$ws.Rows("2:2").Insert()

# Write the updated data values (rows 2-6) for the English word
# scaling analysis re-run.
$rowData = @(26.283883051840483,6.534873231386397,-13.335427761014689,6.304287481517972,422.4126290567733,14.321795995515986,485.86281340968475,9.710804157837853,0.08117316775624991,0.08117316775624991,7.795386774242281,1.1136266820346115,0.776599868600355,-0.5711456461283003)
for ($i = 0; $i -lt $rowData.Length; $i++) { $ws.Cells.Item(2, $i + 1).Value = $rowData[$i] }

$rowData = @(-43.56739527156814,14.435612919278318,-55.74272708316374,15.83946877880295,622.5667789651081,35.20362304295403,482.7188583531113,14.391347530491375,0.0020651019710883796,0.0020651019710883796,23.375762491688974,3.3393946416698532,-0.9943942232114882,-0.3863938875632492)
for ($i = 0; $i -lt $rowData.Length; $i++) { $ws.Cells.Item(3, $i + 1).Value = $rowData[$i] }

$rowData = @(149.9999968377895,13.872931072598096,-26.25539994910415,6.985434185002128,212.82725717686984,27.870525121311985,584.0054502083184,8.363658997566155,1.6243524809476624,1.6243524809476624,23.84705629091389,3.406722327273413,0.9814712855169814,-0.4231625917339956)
for ($i = 0; $i -lt $rowData.Length; $i++) { $ws.Cells.Item(4, $i + 1).Value = $rowData[$i] }

$rowData = @(149.9999931192523,22.73654534990924,-17.420495880944472,13.100204598170073,242.4987023205263,47.088423233717656,535.613158069118,16.180735797511957,0.2329781864739021,0.2329781864739021,28.931134769708052,4.133019252815436,0.9526316899908831,-0.28679405289330295)
for ($i = 0; $i -lt $rowData.Length; $i++) { $ws.Cells.Item(5, $i + 1).Value = $rowData[$i] }

$rowData = @(103.58276186940586,7.534902560791815,-2.967820011403331,7.161744324032855,317.7193284231655,17.96576413844692,538.1133102609925,8.524324585938162,0.2831761154698147,0.2831761154698147,8.841906670537208,1.2631295243624583,0.9663035114294074,-0.30584135700351567)
for ($i = 0; $i -lt $rowData.Length; $i++) { $ws.Cells.Item(6, $i + 1).Value = $rowData[$i] }
